$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2823
$ws.Range("I43").Value = 2430
$ws.Range("J43").Value = 4002
$ws.Range("K43").Value = 2430
$ws.Range("L43").Value = 4002
$ws.Range("M43").Value = -2361
$ws.Range("N43").Value = -4140

$ws.Range("H116").Value = 8332.666999999999
$ws.Range("I116").Value = 8332.666999999999
$ws.Range("K116").Value = 8332.666999999999
$ws.Range("M116").Value = -4890.666999999999

$ws.Range("H132").Value = 4049.2917
$ws.Range("I132").Value = 4200.452
$ws.Range("K132").Value = 12601.356
$ws.Range("M132").Value = -10071.356

$ws.Range("H133").Value = 105911.25
$ws.Range("J133").Value = 105911.25
$ws.Range("L133").Value = 105911.25
$ws.Range("N133").Value = -116031.25

$ws.Range("H138").Value = 3269.4666
$ws.Range("I138").Value = 3498.182
$ws.Range("J138").Value = 3050.6956
$ws.Range("K138").Value = 10494.546
$ws.Range("L138").Value = 9152.086800000001
$ws.Range("M138").Value = -5354.545999999998
$ws.Range("N138").Value = -19432.0868

$ws.Range("H141").Value = 4373.3
$ws.Range("I141").Value = 4389.143
$ws.Range("K141").Value = 13167.429
$ws.Range("M141").Value = -7987.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5136.8237
$ws.Range("I32").Value = 4774
$ws.Range("J32").Value = 6830
$ws.Range("K32").Value = 4774
$ws.Range("L32").Value = 6830
$ws.Range("M32").Value = -4487
$ws.Range("N32").Value = -7404

$ws.Range("H45").Value = 50342.668
$ws.Range("I45").Value = 85805.60000000001
$ws.Range("J45").Value = 6014
$ws.Range("K45").Value = 85805.60000000001
$ws.Range("L45").Value = 6014
$ws.Range("M45").Value = -85428.60000000001
$ws.Range("N45").Value = -6768

$ws.Range("H63").Value = 4065.3333
$ws.Range("I63").Value = 2512.5715
$ws.Range("K63").Value = 2512.5715
$ws.Range("M63").Value = -1826.5715

$ws.Range("H66").Value = 4065.3333
$ws.Range("I66").Value = 2512.5715
$ws.Range("K66").Value = 12562.8575
$ws.Range("M66").Value = -9130.8575

$ws.Range("H132").Value = 3394.8462
$ws.Range("I132").Value = 2033.75
$ws.Range("J132").Value = 3999.7778
$ws.Range("K132").Value = 6101.25
$ws.Range("L132").Value = 11999.3334
$ws.Range("M132").Value = -3571.25
$ws.Range("N132").Value = -17059.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 41685.434
$ws.Range("J138").Value = 41685.434
$ws.Range("L138").Value = 41685.434
$ws.Range("N138").Value = -51965.434

$ws.Range("H139").Value = 51567.918
$ws.Range("J139").Value = 51567.918
$ws.Range("L139").Value = 51567.918
$ws.Range("N139").Value = -61847.918

$ws.Range("H141").Value = 35940.6
$ws.Range("J141").Value = 37425.75
$ws.Range("L141").Value = 37425.75
$ws.Range("N141").Value = -47785.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5357.174
$ws.Range("J31").Value = 7431.1816
$ws.Range("L31").Value = 7431.1816
$ws.Range("N31").Value = -8021.1816

$ws.Range("H34").Value = 5357.174
$ws.Range("J34").Value = 7431.1816
$ws.Range("L34").Value = 7431.1816
$ws.Range("N34").Value = -7835.1816

$ws.Range("H68").Value = 84473.25
$ws.Range("J68").Value = 84473.25
$ws.Range("L68").Value = 84473.25
$ws.Range("N68").Value = -85971.25

$ws.Range("H69").Value = 11250
$ws.Range("I69").Value = 9500
$ws.Range("J69").Value = 13000
$ws.Range("K69").Value = 9500
$ws.Range("L69").Value = 13000
$ws.Range("M69").Value = -8751
$ws.Range("N69").Value = -14498

$ws.Range("H71").Value = 84473.25
$ws.Range("J71").Value = 84473.25
$ws.Range("L71").Value = 253419.75
$ws.Range("N71").Value = -260907.75

$ws.Range("H72").Value = 11250
$ws.Range("I72").Value = 9500
$ws.Range("J72").Value = 13000
$ws.Range("K72").Value = 28500
$ws.Range("L72").Value = 39000
$ws.Range("M72").Value = -24756
$ws.Range("N72").Value = -46488

$ws.Range("H74").Value = 42579.332
$ws.Range("J74").Value = 42579.332
$ws.Range("L74").Value = 42579.332
$ws.Range("N74").Value = -44327.332

$ws.Range("H77").Value = 42579.332
$ws.Range("J77").Value = 42579.332
$ws.Range("L77").Value = 127737.996
$ws.Range("N77").Value = -136473.996

$ws.Range("H104").Value = 49985
$ws.Range("J104").Value = 49985
$ws.Range("L104").Value = 49985
$ws.Range("N104").Value = -55227

$ws.Range("H105").Value = 1235.6
$ws.Range("J105").Value = 2493
$ws.Range("L105").Value = 2493
$ws.Range("N105").Value = -5987

$ws.Range("H134").Value = 2420.15
$ws.Range("I134").Value = 1985.5
$ws.Range("K134").Value = 5956.5
$ws.Range("M134").Value = -3421.5

$ws.Range("H139").Value = 64017
$ws.Range("J139").Value = 64017
$ws.Range("L139").Value = 64017
$ws.Range("N139").Value = -74297

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 62.058823
$ws.Range("I40").Value = 62.058823
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 248.235292
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -179.235292
$ws.Range("N40").ClearContents()

$ws.Range("H95").Value = 7748
$ws.Range("J95").Value = 7748
$ws.Range("L95").Value = 23244
$ws.Range("N95").Value = -27362

$ws.Range("H132").Value = 2473.5557
$ws.Range("I132").Value = 1421.6666
$ws.Range("K132").Value = 12794.9994
$ws.Range("M132").Value = -10264.9994

$ws.Range("H133").Value = 5333
$ws.Range("I133").Value = 2999.5
$ws.Range("J133").Value = 10000
$ws.Range("K133").Value = 8998.5
$ws.Range("L133").Value = 30000
$ws.Range("M133").Value = -3938.5
$ws.Range("N133").Value = -40120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 601.9286
$ws.Range("I2").Value = 855
$ws.Range("J2").Value = 146.4
$ws.Range("K2").Value = 855
$ws.Range("L2").Value = 146.4
$ws.Range("M2").Value = -742
$ws.Range("N2").Value = -372.4

$ws.Range("H122").Value = 5134.08
$ws.Range("I122").Value = 4346
$ws.Range("K122").Value = 13038
$ws.Range("M122").Value = -10588

$ws.Range("H141").Value = 66195
$ws.Range("J141").Value = 72000
$ws.Range("L141").Value = 72000
$ws.Range("N141").Value = -82360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H112").Value = 58787
$ws.Range("J112").Value = 58787
$ws.Range("L112").Value = 58787
$ws.Range("N112").Value = -61741

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H120").Value = 250000
$ws.Range("J120").Value = 250000
$ws.Range("L120").Value = 250000
$ws.Range("N120").Value = -259676

$ws.Range("H136").Value = 3595.28
$ws.Range("I136").Value = 3494.25
$ws.Range("K136").Value = 10482.75
$ws.Range("M136").Value = -7932.75

$ws.Range("H140").Value = 69475.664
$ws.Range("J140").Value = 69475.664
$ws.Range("L140").Value = 69475.664
$ws.Range("N140").Value = -79835.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4249.75
$ws.Range("I81").Value = 4500
$ws.Range("J81").Value = 3999.5
$ws.Range("K81").Value = 9000
$ws.Range("L81").Value = 7999
$ws.Range("M81").Value = -7939
$ws.Range("N81").Value = -10121

$ws.Range("H84").Value = 4249.75
$ws.Range("I84").Value = 4500
$ws.Range("J84").Value = 3999.5
$ws.Range("K84").Value = 45000
$ws.Range("L84").Value = 39995
$ws.Range("M84").Value = -39696
$ws.Range("N84").Value = -50603

$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140

$ws.Range("H140").Value = 100654.11
$ws.Range("J140").Value = 96937.125
$ws.Range("L140").Value = 96937.125
$ws.Range("N140").Value = -107297.125

$ws.Range("H141").Value = 69149.09
$ws.Range("J141").Value = 69999
$ws.Range("L141").Value = 69999
$ws.Range("N141").Value = -80359
